$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) / Volume(1h) (E) figures for the latest crypto snapshot.
# Force text format on just the touched cells so values like "280.19" or
# "5.91%" are stored as literal text, matching the existing sheet layout.
$targets = @("D2","E2","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","E27","D40","E40","D41","E41","D42","E42","E43","E44","D45","E45","E46","E47","E49","E50")
foreach ($addr in $targets) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "280.19"
$ws.Range("E2").Value = "5.91%"
$ws.Range("E3").Value = "1.77%"
$ws.Range("D4").Value = "4.930"
$ws.Range("E4").Value = "4.84%"
$ws.Range("D5").Value = "0.06391"
$ws.Range("E5").Value = "4.93%"
$ws.Range("D6").Value = "6.970"
$ws.Range("E6").Value = "3.57%"
$ws.Range("D7").Value = "3.352"
$ws.Range("E7").Value = "5.76%"
$ws.Range("D8").Value = "0.8848"
$ws.Range("E8").Value = "3.93%"
$ws.Range("D9").Value = "0.9548"
$ws.Range("E9").Value = "4.97%"
$ws.Range("D10").Value = "0.1483"
$ws.Range("E10").Value = "5.22%"
$ws.Range("D11").Value = "0.05191"
$ws.Range("E11").Value = "3.86%"
$ws.Range("D12").Value = "0.07437"
$ws.Range("E12").Value = "4.54%"
$ws.Range("D13").Value = "0.03117"
$ws.Range("E13").Value = "-0.91%"
$ws.Range("E14").Value = "0.31%"
$ws.Range("D15").Value = "0.001560"
$ws.Range("E15").Value = "1.57%"
$ws.Range("D16").Value = "0.0006289"
$ws.Range("E16").Value = "3.97%"
$ws.Range("D17").Value = "0.005875"
$ws.Range("E17").Value = "-2.87%"
$ws.Range("D18").Value = "3.504"
$ws.Range("E18").Value = "1.59%"
$ws.Range("D19").Value = "2.299"
$ws.Range("E19").Value = "5.66%"
$ws.Range("D20").Value = "0.3093"
$ws.Range("D21").Value = "0.1288"
$ws.Range("E21").Value = "0.47%"
$ws.Range("D22").Value = "3.934"
$ws.Range("E22").Value = "-4.34%"
$ws.Range("D23").Value = "0.04326"
$ws.Range("E23").Value = "2.13%"
$ws.Range("D24").Value = "0.001176"
$ws.Range("E24").Value = "-0.22%"
$ws.Range("D25").Value = "0.003665"
$ws.Range("E25").Value = "-9.68%"
$ws.Range("E26").Value = "-0.33%"
$ws.Range("E27").Value = "0.49%"
$ws.Range("D40").Value = "0.04078"
$ws.Range("E40").Value = "3.84%"
$ws.Range("D41").Value = "0.006634"
$ws.Range("E41").Value = "58.52%"
$ws.Range("D42").Value = "0.1176"
$ws.Range("E42").Value = "5.44%"
$ws.Range("E43").Value = "11.78%"
$ws.Range("E44").Value = "7.26%"
$ws.Range("D45").Value = "0.00005245"
$ws.Range("E45").Value = "2.90%"
$ws.Range("E46").Value = "-0.13%"
$ws.Range("E47").Value = "814.01%"
$ws.Range("E49").Value = "-0.13%"
$ws.Range("E50").Value = "-0.20%"
